$d = $word.ActiveDocument

# Locate the paragraph discussing the dispersion of prices (unique text).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "dispers") {
        $target = $p
        break
    }
}

# 1. Reword: "dispersão elevada entre os dados" -> "dispersão entre os
#    preços dos apartamentos".
$target.Range.Find.Execute(
    "dispersão elevada entre os dados", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "dispersão entre os preços dos apartamentos", 2
) | Out-Null

# 2. Merge the following (blank-looking) paragraph into this one by
#    deleting the paragraph mark between them. That following paragraph
#    carries its own pPr (pBdr/spacing/ind/jc) which disappears once the
#    two paragraphs become one, while its empty runs are kept as part of
#    this paragraph.
$nextP = $target.Next()
$mark = $d.Range($target.Range.End - 1, $nextP.Range.End - 1)
$mark.Delete()
